$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows (2-5) hold daily price records for "Albahaca" that get
# reshuffled/reordered (e.g. sorted by date) as part of the weekly
# consolidation. Columns A, B, C, E, F, G, H, I, R stay identical across
# rows, so only D (Fecha), J (Volumen), K (Precio minimo), L (Precio
# maximo), M (Precio promedio ponderado), N (Unidad de comercializacion),
# O (Origen), P (Precio $/Kg) and Q (Kg o Unidades) need to move.

$ws.Range("D2").Value = 44687
$ws.Range("J2").Value = 160
$ws.Range("K2").Value = 3000
$ws.Range("L2").Value = 3500
$ws.Range("M2").Value = 3250
$ws.Range("N2").Value = "$/docena de matas"
$ws.Range("O2").Value = "Región Metropolitana"
$ws.Range("P2").Value = 542
$ws.Range("Q2").Value = 6

$ws.Range("D3").Value = 44221
$ws.Range("J3").Value = 250
$ws.Range("K3").Value = 1300
$ws.Range("L3").Value = 1500
$ws.Range("M3").Value = 1420
$ws.Range("N3").Value = "$/atado"
$ws.Range("O3").Value = "Provincia de Diguillín"
$ws.Range("P3").Value = 1420
$ws.Range("Q3").Value = 1

$ws.Range("D4").Value = 44691
$ws.Range("J4").Value = 100
$ws.Range("K4").Value = 3000
$ws.Range("L4").Value = 3500
$ws.Range("M4").Value = 3250
$ws.Range("N4").Value = "$/docena de matas"
$ws.Range("O4").Value = "Región Metropolitana"
$ws.Range("P4").Value = 542
$ws.Range("Q4").Value = 6

$ws.Range("D5").Value = 44692
$ws.Range("J5").Value = 120
$ws.Range("K5").Value = 3000
$ws.Range("L5").Value = 3500
$ws.Range("M5").Value = 3250
$ws.Range("N5").Value = "$/docena de matas"
$ws.Range("O5").Value = "Región Metropolitana"
$ws.Range("P5").Value = 542
$ws.Range("Q5").Value = 6
